$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.285.89"
$ws.Range("E2").Value = "  -1.85%  "

$ws.Range("D3").Value = "3.067.58"
$ws.Range("E3").Value = "  -2.73%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +7.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "615.36"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.70%  "

$ws.Range("E7").Value = "  -6.90%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.358"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.00%  "

$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").Value = "3.065.81"
$ws.Range("E10").Value = "  -2.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.705"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -6.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.197"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.88"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.86%  "

$ws.Range("D15").Value = "89.305.04"
$ws.Range("E15").Value = "  -1.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.33"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -6.66%  "

$ws.Range("D17").Value = "3.637.38"
$ws.Range("E17").Value = "  -2.70%  "

$ws.Range("D18").Value = "3.087.09"
$ws.Range("E18").Value = "  -1.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.74"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.46%  "

$ws.Range("E20").Value = "  -0.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.67"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "428.88"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -9.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.35"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.67"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.55"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -6.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "86.41"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -9.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.62"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -6.06%  "

$ws.Range("D28").Value = "3.244.61"
$ws.Range("E28").Value = "  -2.31%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.11"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +10.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.96"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.155"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.75%  "

$ws.Range("E33").Value = "  -14.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.44"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -6.58%  "

$ws.Range("E35").Value = "  +2.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.02"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "490.09"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.97%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.60"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.80%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.24"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -6.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.64"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +52.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0890"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.04"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.88%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.395"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -8.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "152.95"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.95%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.83"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -7.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.666"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -8.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.38"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("E51").Value = "  -5.77%  "
